$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add missing dw_filled_tube_food_control_mass (column E) values for rows 14-16
$ws.Range("E14").Value = 5733.2
$ws.Range("E15").Value = 5745
$ws.Range("E16").Value = 5695.1

# Update the selected cell to reflect the saved cursor position
$ws.Range("E14").Select()
